# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 14:59"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 6461353
$ws.Range("C4").Value = 1103
$ws.Range("E4").Value = 2542001

# --- Row 13: Argentina ---
$ws.Range("D13").Value = 357388
$ws.Range("E13").Value = 111492
$ws.Range("G13").Value = 53
$ws.Range("H13").Value = 9912

# --- Row 19: Arabia Saudita ---
$ws.Range("B19").Value = 321456
$ws.Range("C19").Value = 768
$ws.Range("D19").Value = 297623
$ws.Range("E19").Value = 19726
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 4107

# --- Row 24: Alemania ---
$ws.Range("B24").Value = 252069
$ws.Range("C24").Value = 345
$ws.Range("E24").Value = 15668

# --- Rows 41/42: Suecia overtakes China, rows swap (sorted desc by Casos totales) ---
# Row 41 becomes Suecia with updated figures
$ws.Range("A41").Value = "Suecia"
$ws.Range("B41").Value = 85558
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 5
$ws.Range("H41").Value = 5837

# Row 42 becomes China with its (unchanged) figures
$ws.Range("A42").Value = "China"
$ws.Range("B42").Value = 85134
$ws.Range("C42").Value = 12
$ws.Range("D42").Value = 80320
$ws.Range("E42").Value = 180
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 4634

# --- Row 46: Bielorrusia ---
$ws.Range("B46").Value = 73031
$ws.Range("C46").Value = 172
$ws.Range("D46").Value = 71883
$ws.Range("E46").Value = 432
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 716

# --- Row 73: Estado de Palestina ---
$ws.Range("E73").Value = 9325
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 184

# --- Row 83: Dinamarca ---
$ws.Range("B83").Value = 18113
$ws.Range("C83").Value = 230
$ws.Range("D83").Value = 15833
$ws.Range("E83").Value = 1652
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 628

# --- Row 87: Senegal ---
$ws.Range("B87").Value = 14014
$ws.Range("C87").Value = 27
$ws.Range("D87").Value = 10037
$ws.Range("E87").Value = 3686
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 291
